$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 5 (duplicate fastq file entry), shifting rows 6:35 up to 5:34.
$ws.Rows.Item(5).Delete()

# Reflect the resulting selection state, matching Excel's default behavior after
# deleting a row via the row header (whole row 5 selected post-delete).
$ws.Range("A5:XFD5").Select()
